# Adapt the "flash分配" address table from absolute flash addresses (0x2Fxx)
# to an offset + sector-base-address scheme, matching the wifi-bt warm/cold
# branch's log layout:
#   - column E header "地址" -> "偏移地址", values become small offsets
#   - new column F "扇区基地址" with the sector base address macro name

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row ---------------------------------------------------------
$ws.Range("E1").Value = "偏移地址"
$ws.Range("F1").Value = "扇区基地址"

# --- Column E: absolute address -> offset-only address -----------------
$ws.Range("E2").Value  = "0x00~0x01"
$ws.Range("E3").Value  = "0x02"
$ws.Range("E4").Value  = "0x03~0x04"
$ws.Range("E5").Value  = "0x05"
# E6 stays blank (merged into E5:E6)
$ws.Range("E7").Value  = "0x06"
$ws.Range("E8").Value  = "0x07"
$ws.Range("E9").Value  = "0x08"
$ws.Range("E10").Value = "0x09"
$ws.Range("E11").Value = "0x0A"
$ws.Range("E12").Value = "0x0B"
$ws.Range("E13").Value = "0x00"
$ws.Range("E14").Value = "0x01"

# --- Column F: new "sector base address" column -------------------------
$ws.Range("F2:F12").Value = "USER_PARAMETER_START_SECTOR_ADDRESS0"
$ws.Range("F13:F14").Value = "USER_PARAMETER_START_SECTOR_ADDRESS1"

# --- Formatting: give column F the same bordered / left-aligned style
#     already used throughout the table (copy from column B). -----------
$ws.Range("B1:B14").Copy()
$ws.Range("F1:F14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# E5/E6 (merged) now share one unified bordered style, same as before but
# consolidated - copy E5's (already bordered) style onto E6 so both halves
# of the merged cell match.
$ws.Range("E5").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column F width, matching the new column's content width -----------
$ws.Columns.Item(6).ColumnWidth = 38.14

# --- Selection ends on the last edited cell, like the source commit -----
$ws.Range("F14").Select()
